$updates = @(
    @{ Cell = "D2"; Value = "56.580.80" }
    @{ Cell = "E2"; Value = "  +4.10%  " }
    @{ Cell = "D3"; Value = "2.482.40" }
    @{ Cell = "E3"; Value = "  +0.53%  " }
    @{ Cell = "E4"; Value = "  -0.07%  " }
    @{ Cell = "D5"; Value = "487.12" }
    @{ Cell = "E5"; Value = "  +4.05%  " }
    @{ Cell = "D6"; Value = "145.61" }
    @{ Cell = "E6"; Value = "  +8.95%  " }
    @{ Cell = "D7"; Value = "0.998" }
    @{ Cell = "E7"; Value = "  +0.02%  " }
    @{ Cell = "D8"; Value = "0.512" }
    @{ Cell = "E8"; Value = "  +4.15%  " }
    @{ Cell = "D9"; Value = "2.501.38" }
    @{ Cell = "E9"; Value = "  +1.31%  " }
    @{ Cell = "D10"; Value = "5.77" }
    @{ Cell = "E10"; Value = "  +7.95%  " }
    @{ Cell = "D11"; Value = "0.0971" }
    @{ Cell = "E11"; Value = "  +0.87%  " }
    @{ Cell = "D12"; Value = "0.331" }
    @{ Cell = "E12"; Value = "  +3.70%  " }
    @{ Cell = "D13"; Value = "0.123" }
    @{ Cell = "E13"; Value = "  +0.92%  " }
    @{ Cell = "D14"; Value = "2.916.39" }
    @{ Cell = "E14"; Value = "  +0.79%  " }
    @{ Cell = "D15"; Value = "56.335.75" }
    @{ Cell = "E15"; Value = "  +3.83%  " }
    @{ Cell = "D16"; Value = "21.19" }
    @{ Cell = "E16"; Value = "  +6.14%  " }
    @{ Cell = "D17"; Value = "0.0000136" }
    @{ Cell = "E17"; Value = "  +2.96%  " }
    @{ Cell = "D18"; Value = "2.501.44" }
    @{ Cell = "E18"; Value = "  +1.29%  " }
    @{ Cell = "D19"; Value = "4.54" }
    @{ Cell = "E19"; Value = "  +7.07%  " }
    @{ Cell = "D20"; Value = "10.17" }
    @{ Cell = "E20"; Value = "  +6.54%  " }
    @{ Cell = "D21"; Value = "319.81" }
    @{ Cell = "E21"; Value = "  +3.03%  " }
    @{ Cell = "D22"; Value = "0.996" }
    @{ Cell = "E22"; Value = "  -0.62%  " }
    @{ Cell = "D23"; Value = "5.82" }
    @{ Cell = "E23"; Value = "  +8.19%  " }
    @{ Cell = "D24"; Value = "58.54" }
    @{ Cell = "E24"; Value = "  +2.82%  " }
    @{ Cell = "D25"; Value = "0.411" }
    @{ Cell = "E25"; Value = "  +5.96%  " }
    @{ Cell = "D26"; Value = "0.165" }
    @{ Cell = "E26"; Value = "  +7.16%  " }
    @{ Cell = "E27"; Value = "  -1.04%  " }
    @{ Cell = "D28"; Value = "2.611.97" }
    @{ Cell = "E28"; Value = "  +2.47%  " }
    @{ Cell = "D29"; Value = "7.68" }
    @{ Cell = "E29"; Value = "  +5.54%  " }
    @{ Cell = "D30"; Value = "0.0₃0790" }
    @{ Cell = "E30"; Value = "  +8.21%  " }
    @{ Cell = "E31"; Value = "  +0.14%  " }
    @{ Cell = "D32"; Value = "149.08" }
    @{ Cell = "E32"; Value = "  -1.00%  " }
    @{ Cell = "D33"; Value = "18.24" }
    @{ Cell = "E33"; Value = "  +1.81%  " }
    @{ Cell = "D34"; Value = "1.50" }
    @{ Cell = "E34"; Value = "  +4.01%  " }
    @{ Cell = "D35"; Value = "5.20" }
    @{ Cell = "E35"; Value = "  +2.79%  " }
    @{ Cell = "E36"; Value = "  +6.80%  " }
    @{ Cell = "D37"; Value = "3.72" }
    @{ Cell = "E37"; Value = "  +4.03%  " }
    @{ Cell = "D38"; Value = "0.862" }
    @{ Cell = "E38"; Value = "  +6.08%  " }
    @{ Cell = "D39"; Value = "34.18" }
    @{ Cell = "E39"; Value = "  +1.30%  " }
    @{ Cell = "D40"; Value = "3.52" }
    @{ Cell = "E40"; Value = "  +6.54%  " }
    @{ Cell = "D41"; Value = "0.614" }
    @{ Cell = "E41"; Value = "  +1.54%  " }
    @{ Cell = "D42"; Value = "0.0557" }
    @{ Cell = "E42"; Value = "  +4.64%  " }
    @{ Cell = "D43"; Value = "0.997" }
    @{ Cell = "E43"; Value = "  +0.12%  " }
    @{ Cell = "D44"; Value = "1.33" }
    @{ Cell = "E44"; Value = "  +6.49%  " }
    @{ Cell = "D45"; Value = "4.85" }
    @{ Cell = "E45"; Value = "  +12.49%  " }
    @{ Cell = "D46"; Value = "259.47" }
    @{ Cell = "E46"; Value = "  +15.08%  " }
    @{ Cell = "B47"; Value = "WhiteBITCoin" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt" }
    @{ Cell = "D47"; Value = "10.18" }
    @{ Cell = "E47"; Value = "  -0.14%  " }
    @{ Cell = "B48"; Value = "VeChain" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "D48"; Value = "0.0228" }
    @{ Cell = "E48"; Value = "  +3.83%  " }
    @{ Cell = "D49"; Value = "0.0911" }
    @{ Cell = "E49"; Value = "  +3.63%  " }
    @{ Cell = "D50"; Value = "1.907.53" }
    @{ Cell = "E50"; Value = "  -2.53%  " }
    @{ Cell = "D51"; Value = "17.61" }
    @{ Cell = "E51"; Value = "  +5.06%  " }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}
